# Rewrites "ODI Batting Extra" (sheet 4) rows 2-70 with the extra-scrape
# batting/bowling data (MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6,
# PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH). Text-looking numeric values are
# entered with a leading apostrophe so Excel keeps them as text instead
# of coercing to Number/Percent; ClearFormats() strips the resulting
# "quote prefix" style flag so formatting stays untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Row 2
$ws.Cells.Item(2,1).Value = "'3129"
$ws.Cells.Item(2,1).ClearFormats()
$ws.Cells.Item(2,2).Value = 10
$ws.Cells.Item(2,2).ClearFormats()
$ws.Cells.Item(2,3).Value = "'"
$ws.Cells.Item(2,3).ClearFormats()
$ws.Cells.Item(2,4).Value = "'"
$ws.Cells.Item(2,4).ClearFormats()
$ws.Cells.Item(2,5).Value = "'"
$ws.Cells.Item(2,5).ClearFormats()
$ws.Cells.Item(2,6).Value = "'NO"
$ws.Cells.Item(2,6).ClearFormats()

# Row 3
$ws.Cells.Item(3,1).Value = "'3557"
$ws.Cells.Item(3,1).ClearFormats()
$ws.Cells.Item(3,2).Value = "'"
$ws.Cells.Item(3,2).ClearFormats()
$ws.Cells.Item(3,3).Value = "'"
$ws.Cells.Item(3,3).ClearFormats()
$ws.Cells.Item(3,4).Value = "'"
$ws.Cells.Item(3,4).ClearFormats()
$ws.Cells.Item(3,5).Value = "'"
$ws.Cells.Item(3,5).ClearFormats()
$ws.Cells.Item(3,6).Value = "'NO"
$ws.Cells.Item(3,6).ClearFormats()

# Row 4
$ws.Cells.Item(4,1).Value = "'3694"
$ws.Cells.Item(4,1).ClearFormats()
$ws.Cells.Item(4,2).Value = "'"
$ws.Cells.Item(4,2).ClearFormats()
$ws.Cells.Item(4,3).Value = "'"
$ws.Cells.Item(4,3).ClearFormats()
$ws.Cells.Item(4,4).Value = "'"
$ws.Cells.Item(4,4).ClearFormats()
$ws.Cells.Item(4,5).Value = "'"
$ws.Cells.Item(4,5).ClearFormats()
$ws.Cells.Item(4,6).Value = "'NO"
$ws.Cells.Item(4,6).ClearFormats()

# Row 5
$ws.Cells.Item(5,1).Value = "'3695"
$ws.Cells.Item(5,1).ClearFormats()
$ws.Cells.Item(5,2).Value = 11
$ws.Cells.Item(5,2).ClearFormats()
$ws.Cells.Item(5,3).Value = "'0"
$ws.Cells.Item(5,3).ClearFormats()
$ws.Cells.Item(5,4).Value = "'0"
$ws.Cells.Item(5,4).ClearFormats()
$ws.Cells.Item(5,5).Value = "'"
$ws.Cells.Item(5,5).ClearFormats()
$ws.Cells.Item(5,6).Value = "'NO"
$ws.Cells.Item(5,6).ClearFormats()

# Row 6
$ws.Cells.Item(6,1).Value = "'3697"
$ws.Cells.Item(6,1).ClearFormats()
$ws.Cells.Item(6,2).Value = 11
$ws.Cells.Item(6,2).ClearFormats()
$ws.Cells.Item(6,3).Value = "'"
$ws.Cells.Item(6,3).ClearFormats()
$ws.Cells.Item(6,4).Value = "'"
$ws.Cells.Item(6,4).ClearFormats()
$ws.Cells.Item(6,5).Value = "'"
$ws.Cells.Item(6,5).ClearFormats()
$ws.Cells.Item(6,6).Value = "'NO"
$ws.Cells.Item(6,6).ClearFormats()

# Row 7
$ws.Cells.Item(7,1).Value = "'3700"
$ws.Cells.Item(7,1).ClearFormats()
$ws.Cells.Item(7,2).Value = "'"
$ws.Cells.Item(7,2).ClearFormats()
$ws.Cells.Item(7,3).Value = "'"
$ws.Cells.Item(7,3).ClearFormats()
$ws.Cells.Item(7,4).Value = "'"
$ws.Cells.Item(7,4).ClearFormats()
$ws.Cells.Item(7,5).Value = "'"
$ws.Cells.Item(7,5).ClearFormats()
$ws.Cells.Item(7,6).Value = "'NO"
$ws.Cells.Item(7,6).ClearFormats()

# Row 8
$ws.Cells.Item(8,1).Value = "'3741"
$ws.Cells.Item(8,1).ClearFormats()
$ws.Cells.Item(8,2).Value = "'"
$ws.Cells.Item(8,2).ClearFormats()
$ws.Cells.Item(8,3).Value = "'"
$ws.Cells.Item(8,3).ClearFormats()
$ws.Cells.Item(8,4).Value = "'"
$ws.Cells.Item(8,4).ClearFormats()
$ws.Cells.Item(8,5).Value = "'"
$ws.Cells.Item(8,5).ClearFormats()
$ws.Cells.Item(8,6).Value = "'NO"
$ws.Cells.Item(8,6).ClearFormats()

# Row 9
$ws.Cells.Item(9,1).Value = "'3746"
$ws.Cells.Item(9,1).ClearFormats()
$ws.Cells.Item(9,2).Value = "'"
$ws.Cells.Item(9,2).ClearFormats()
$ws.Cells.Item(9,3).Value = "'"
$ws.Cells.Item(9,3).ClearFormats()
$ws.Cells.Item(9,4).Value = "'"
$ws.Cells.Item(9,4).ClearFormats()
$ws.Cells.Item(9,5).Value = "'"
$ws.Cells.Item(9,5).ClearFormats()
$ws.Cells.Item(9,6).Value = "'NO"
$ws.Cells.Item(9,6).ClearFormats()

# Row 10
$ws.Cells.Item(10,1).Value = "'3749"
$ws.Cells.Item(10,1).ClearFormats()
$ws.Cells.Item(10,2).Value = 11
$ws.Cells.Item(10,2).ClearFormats()
$ws.Cells.Item(10,3).Value = "'"
$ws.Cells.Item(10,3).ClearFormats()
$ws.Cells.Item(10,4).Value = "'"
$ws.Cells.Item(10,4).ClearFormats()
$ws.Cells.Item(10,5).Value = "'"
$ws.Cells.Item(10,5).ClearFormats()
$ws.Cells.Item(10,6).Value = "'NO"
$ws.Cells.Item(10,6).ClearFormats()

# Row 11
$ws.Cells.Item(11,1).Value = "'3773"
$ws.Cells.Item(11,1).ClearFormats()
$ws.Cells.Item(11,2).Value = 11
$ws.Cells.Item(11,2).ClearFormats()
$ws.Cells.Item(11,3).Value = "'"
$ws.Cells.Item(11,3).ClearFormats()
$ws.Cells.Item(11,4).Value = "'"
$ws.Cells.Item(11,4).ClearFormats()
$ws.Cells.Item(11,5).Value = "'"
$ws.Cells.Item(11,5).ClearFormats()
$ws.Cells.Item(11,6).Value = "'NO"
$ws.Cells.Item(11,6).ClearFormats()

# Row 12
$ws.Cells.Item(12,1).Value = "'3792"
$ws.Cells.Item(12,1).ClearFormats()
$ws.Cells.Item(12,2).Value = "'"
$ws.Cells.Item(12,2).ClearFormats()
$ws.Cells.Item(12,3).Value = "'"
$ws.Cells.Item(12,3).ClearFormats()
$ws.Cells.Item(12,4).Value = "'"
$ws.Cells.Item(12,4).ClearFormats()
$ws.Cells.Item(12,5).Value = "'"
$ws.Cells.Item(12,5).ClearFormats()
$ws.Cells.Item(12,6).Value = "'NO"
$ws.Cells.Item(12,6).ClearFormats()

# Row 13
$ws.Cells.Item(13,1).Value = "'3795"
$ws.Cells.Item(13,1).ClearFormats()
$ws.Cells.Item(13,2).Value = 11
$ws.Cells.Item(13,2).ClearFormats()
$ws.Cells.Item(13,3).Value = "'"
$ws.Cells.Item(13,3).ClearFormats()
$ws.Cells.Item(13,4).Value = "'"
$ws.Cells.Item(13,4).ClearFormats()
$ws.Cells.Item(13,5).Value = "'"
$ws.Cells.Item(13,5).ClearFormats()
$ws.Cells.Item(13,6).Value = "'NO"
$ws.Cells.Item(13,6).ClearFormats()

# Row 14
$ws.Cells.Item(14,1).Value = "'3796"
$ws.Cells.Item(14,1).ClearFormats()
$ws.Cells.Item(14,2).Value = 10
$ws.Cells.Item(14,2).ClearFormats()
$ws.Cells.Item(14,3).Value = "'"
$ws.Cells.Item(14,3).ClearFormats()
$ws.Cells.Item(14,4).Value = "'"
$ws.Cells.Item(14,4).ClearFormats()
$ws.Cells.Item(14,5).Value = "'"
$ws.Cells.Item(14,5).ClearFormats()
$ws.Cells.Item(14,6).Value = "'NO"
$ws.Cells.Item(14,6).ClearFormats()

# Row 15
$ws.Cells.Item(15,1).Value = "'3874"
$ws.Cells.Item(15,1).ClearFormats()
$ws.Cells.Item(15,2).Value = "'"
$ws.Cells.Item(15,2).ClearFormats()
$ws.Cells.Item(15,3).Value = "'"
$ws.Cells.Item(15,3).ClearFormats()
$ws.Cells.Item(15,4).Value = "'"
$ws.Cells.Item(15,4).ClearFormats()
$ws.Cells.Item(15,5).Value = "'"
$ws.Cells.Item(15,5).ClearFormats()
$ws.Cells.Item(15,6).Value = "'NO"
$ws.Cells.Item(15,6).ClearFormats()

# Row 16
$ws.Cells.Item(16,1).Value = "'3884"
$ws.Cells.Item(16,1).ClearFormats()
$ws.Cells.Item(16,2).Value = 10
$ws.Cells.Item(16,2).ClearFormats()
$ws.Cells.Item(16,3).Value = "'"
$ws.Cells.Item(16,3).ClearFormats()
$ws.Cells.Item(16,4).Value = "'"
$ws.Cells.Item(16,4).ClearFormats()
$ws.Cells.Item(16,5).Value = "'"
$ws.Cells.Item(16,5).ClearFormats()
$ws.Cells.Item(16,6).Value = "'NO"
$ws.Cells.Item(16,6).ClearFormats()

# Row 17
$ws.Cells.Item(17,1).Value = "'3886"
$ws.Cells.Item(17,1).ClearFormats()
$ws.Cells.Item(17,2).Value = 10
$ws.Cells.Item(17,2).ClearFormats()
$ws.Cells.Item(17,3).Value = "'1"
$ws.Cells.Item(17,3).ClearFormats()
$ws.Cells.Item(17,4).Value = "'0"
$ws.Cells.Item(17,4).ClearFormats()
$ws.Cells.Item(17,5).Value = "'2.62%"
$ws.Cells.Item(17,5).ClearFormats()
$ws.Cells.Item(17,6).Value = "'NO"
$ws.Cells.Item(17,6).ClearFormats()

# Row 18
$ws.Cells.Item(18,1).Value = "'3888"
$ws.Cells.Item(18,1).ClearFormats()
$ws.Cells.Item(18,2).Value = "'"
$ws.Cells.Item(18,2).ClearFormats()
$ws.Cells.Item(18,3).Value = "'"
$ws.Cells.Item(18,3).ClearFormats()
$ws.Cells.Item(18,4).Value = "'"
$ws.Cells.Item(18,4).ClearFormats()
$ws.Cells.Item(18,5).Value = "'"
$ws.Cells.Item(18,5).ClearFormats()
$ws.Cells.Item(18,6).Value = "'NO"
$ws.Cells.Item(18,6).ClearFormats()

# Row 19
$ws.Cells.Item(19,1).Value = "'3893"
$ws.Cells.Item(19,1).ClearFormats()
$ws.Cells.Item(19,2).Value = 11
$ws.Cells.Item(19,2).ClearFormats()
$ws.Cells.Item(19,3).Value = "'"
$ws.Cells.Item(19,3).ClearFormats()
$ws.Cells.Item(19,4).Value = "'"
$ws.Cells.Item(19,4).ClearFormats()
$ws.Cells.Item(19,5).Value = "'"
$ws.Cells.Item(19,5).ClearFormats()
$ws.Cells.Item(19,6).Value = "'NO"
$ws.Cells.Item(19,6).ClearFormats()

# Row 20
$ws.Cells.Item(20,1).Value = "'3894"
$ws.Cells.Item(20,1).ClearFormats()
$ws.Cells.Item(20,2).Value = 11
$ws.Cells.Item(20,2).ClearFormats()
$ws.Cells.Item(20,3).Value = "'"
$ws.Cells.Item(20,3).ClearFormats()
$ws.Cells.Item(20,4).Value = "'"
$ws.Cells.Item(20,4).ClearFormats()
$ws.Cells.Item(20,5).Value = "'"
$ws.Cells.Item(20,5).ClearFormats()
$ws.Cells.Item(20,6).Value = "'NO"
$ws.Cells.Item(20,6).ClearFormats()

# Row 21
$ws.Cells.Item(21,1).Value = "'3896"
$ws.Cells.Item(21,1).ClearFormats()
$ws.Cells.Item(21,2).Value = "'"
$ws.Cells.Item(21,2).ClearFormats()
$ws.Cells.Item(21,3).Value = "'"
$ws.Cells.Item(21,3).ClearFormats()
$ws.Cells.Item(21,4).Value = "'"
$ws.Cells.Item(21,4).ClearFormats()
$ws.Cells.Item(21,5).Value = "'"
$ws.Cells.Item(21,5).ClearFormats()
$ws.Cells.Item(21,6).Value = "'NO"
$ws.Cells.Item(21,6).ClearFormats()

# Row 22
$ws.Cells.Item(22,1).Value = "'3898"
$ws.Cells.Item(22,1).ClearFormats()
$ws.Cells.Item(22,2).Value = "'"
$ws.Cells.Item(22,2).ClearFormats()
$ws.Cells.Item(22,3).Value = "'"
$ws.Cells.Item(22,3).ClearFormats()
$ws.Cells.Item(22,4).Value = "'"
$ws.Cells.Item(22,4).ClearFormats()
$ws.Cells.Item(22,5).Value = "'"
$ws.Cells.Item(22,5).ClearFormats()
$ws.Cells.Item(22,6).Value = "'NO"
$ws.Cells.Item(22,6).ClearFormats()

# Row 23
$ws.Cells.Item(23,1).Value = "'3903"
$ws.Cells.Item(23,1).ClearFormats()
$ws.Cells.Item(23,2).Value = 11
$ws.Cells.Item(23,2).ClearFormats()
$ws.Cells.Item(23,3).Value = "'0"
$ws.Cells.Item(23,3).ClearFormats()
$ws.Cells.Item(23,4).Value = "'0"
$ws.Cells.Item(23,4).ClearFormats()
$ws.Cells.Item(23,5).Value = "'"
$ws.Cells.Item(23,5).ClearFormats()
$ws.Cells.Item(23,6).Value = "'NO"
$ws.Cells.Item(23,6).ClearFormats()

# Row 24
$ws.Cells.Item(24,1).Value = "'3905"
$ws.Cells.Item(24,1).ClearFormats()
$ws.Cells.Item(24,2).Value = "'"
$ws.Cells.Item(24,2).ClearFormats()
$ws.Cells.Item(24,3).Value = "'"
$ws.Cells.Item(24,3).ClearFormats()
$ws.Cells.Item(24,4).Value = "'"
$ws.Cells.Item(24,4).ClearFormats()
$ws.Cells.Item(24,5).Value = "'"
$ws.Cells.Item(24,5).ClearFormats()
$ws.Cells.Item(24,6).Value = "'NO"
$ws.Cells.Item(24,6).ClearFormats()

# Row 25
$ws.Cells.Item(25,1).Value = "'3909"
$ws.Cells.Item(25,1).ClearFormats()
$ws.Cells.Item(25,2).Value = "'"
$ws.Cells.Item(25,2).ClearFormats()
$ws.Cells.Item(25,3).Value = "'"
$ws.Cells.Item(25,3).ClearFormats()
$ws.Cells.Item(25,4).Value = "'"
$ws.Cells.Item(25,4).ClearFormats()
$ws.Cells.Item(25,5).Value = "'"
$ws.Cells.Item(25,5).ClearFormats()
$ws.Cells.Item(25,6).Value = "'NO"
$ws.Cells.Item(25,6).ClearFormats()

# Row 26
$ws.Cells.Item(26,1).Value = "'3923"
$ws.Cells.Item(26,1).ClearFormats()
$ws.Cells.Item(26,2).Value = "'"
$ws.Cells.Item(26,2).ClearFormats()
$ws.Cells.Item(26,3).Value = "'"
$ws.Cells.Item(26,3).ClearFormats()
$ws.Cells.Item(26,4).Value = "'"
$ws.Cells.Item(26,4).ClearFormats()
$ws.Cells.Item(26,5).Value = "'"
$ws.Cells.Item(26,5).ClearFormats()
$ws.Cells.Item(26,6).Value = "'NO"
$ws.Cells.Item(26,6).ClearFormats()

# Row 27
$ws.Cells.Item(27,1).Value = "'3927"
$ws.Cells.Item(27,1).ClearFormats()
$ws.Cells.Item(27,2).Value = 11
$ws.Cells.Item(27,2).ClearFormats()
$ws.Cells.Item(27,3).Value = "'"
$ws.Cells.Item(27,3).ClearFormats()
$ws.Cells.Item(27,4).Value = "'"
$ws.Cells.Item(27,4).ClearFormats()
$ws.Cells.Item(27,5).Value = "'"
$ws.Cells.Item(27,5).ClearFormats()
$ws.Cells.Item(27,6).Value = "'NO"
$ws.Cells.Item(27,6).ClearFormats()

# Row 28
$ws.Cells.Item(28,1).Value = "'3966"
$ws.Cells.Item(28,1).ClearFormats()
$ws.Cells.Item(28,2).Value = "'"
$ws.Cells.Item(28,2).ClearFormats()
$ws.Cells.Item(28,3).Value = "'"
$ws.Cells.Item(28,3).ClearFormats()
$ws.Cells.Item(28,4).Value = "'"
$ws.Cells.Item(28,4).ClearFormats()
$ws.Cells.Item(28,5).Value = "'"
$ws.Cells.Item(28,5).ClearFormats()
$ws.Cells.Item(28,6).Value = "'NO"
$ws.Cells.Item(28,6).ClearFormats()

# Row 29
$ws.Cells.Item(29,1).Value = "'3967"
$ws.Cells.Item(29,1).ClearFormats()
$ws.Cells.Item(29,2).Value = "'"
$ws.Cells.Item(29,2).ClearFormats()
$ws.Cells.Item(29,3).Value = "'"
$ws.Cells.Item(29,3).ClearFormats()
$ws.Cells.Item(29,4).Value = "'"
$ws.Cells.Item(29,4).ClearFormats()
$ws.Cells.Item(29,5).Value = "'"
$ws.Cells.Item(29,5).ClearFormats()
$ws.Cells.Item(29,6).Value = "'NO"
$ws.Cells.Item(29,6).ClearFormats()

# Row 30
$ws.Cells.Item(30,1).Value = "'3968"
$ws.Cells.Item(30,1).ClearFormats()
$ws.Cells.Item(30,2).Value = 10
$ws.Cells.Item(30,2).ClearFormats()
$ws.Cells.Item(30,3).Value = "'"
$ws.Cells.Item(30,3).ClearFormats()
$ws.Cells.Item(30,4).Value = "'"
$ws.Cells.Item(30,4).ClearFormats()
$ws.Cells.Item(30,5).Value = "'"
$ws.Cells.Item(30,5).ClearFormats()
$ws.Cells.Item(30,6).Value = "'NO"
$ws.Cells.Item(30,6).ClearFormats()

# Row 31
$ws.Cells.Item(31,1).Value = "'3973"
$ws.Cells.Item(31,1).ClearFormats()
$ws.Cells.Item(31,2).Value = 11
$ws.Cells.Item(31,2).ClearFormats()
$ws.Cells.Item(31,3).Value = "'"
$ws.Cells.Item(31,3).ClearFormats()
$ws.Cells.Item(31,4).Value = "'"
$ws.Cells.Item(31,4).ClearFormats()
$ws.Cells.Item(31,5).Value = "'"
$ws.Cells.Item(31,5).ClearFormats()
$ws.Cells.Item(31,6).Value = "'NO"
$ws.Cells.Item(31,6).ClearFormats()

# Row 32
$ws.Cells.Item(32,1).Value = "'3975"
$ws.Cells.Item(32,1).ClearFormats()
$ws.Cells.Item(32,2).Value = 11
$ws.Cells.Item(32,2).ClearFormats()
$ws.Cells.Item(32,3).Value = "'"
$ws.Cells.Item(32,3).ClearFormats()
$ws.Cells.Item(32,4).Value = "'"
$ws.Cells.Item(32,4).ClearFormats()
$ws.Cells.Item(32,5).Value = "'"
$ws.Cells.Item(32,5).ClearFormats()
$ws.Cells.Item(32,6).Value = "'NO"
$ws.Cells.Item(32,6).ClearFormats()

# Row 33
$ws.Cells.Item(33,1).Value = "'3977"
$ws.Cells.Item(33,1).ClearFormats()
$ws.Cells.Item(33,2).Value = "'"
$ws.Cells.Item(33,2).ClearFormats()
$ws.Cells.Item(33,3).Value = "'"
$ws.Cells.Item(33,3).ClearFormats()
$ws.Cells.Item(33,4).Value = "'"
$ws.Cells.Item(33,4).ClearFormats()
$ws.Cells.Item(33,5).Value = "'"
$ws.Cells.Item(33,5).ClearFormats()
$ws.Cells.Item(33,6).Value = "'NO"
$ws.Cells.Item(33,6).ClearFormats()

# Row 34
$ws.Cells.Item(34,1).Value = "'3981"
$ws.Cells.Item(34,1).ClearFormats()
$ws.Cells.Item(34,2).Value = "'"
$ws.Cells.Item(34,2).ClearFormats()
$ws.Cells.Item(34,3).Value = "'"
$ws.Cells.Item(34,3).ClearFormats()
$ws.Cells.Item(34,4).Value = "'"
$ws.Cells.Item(34,4).ClearFormats()
$ws.Cells.Item(34,5).Value = "'"
$ws.Cells.Item(34,5).ClearFormats()
$ws.Cells.Item(34,6).Value = "'NO"
$ws.Cells.Item(34,6).ClearFormats()

# Row 35
$ws.Cells.Item(35,1).Value = "'3984"
$ws.Cells.Item(35,1).ClearFormats()
$ws.Cells.Item(35,2).Value = 11
$ws.Cells.Item(35,2).ClearFormats()
$ws.Cells.Item(35,3).Value = "'1"
$ws.Cells.Item(35,3).ClearFormats()
$ws.Cells.Item(35,4).Value = "'0"
$ws.Cells.Item(35,4).ClearFormats()
$ws.Cells.Item(35,5).Value = "'4.33%"
$ws.Cells.Item(35,5).ClearFormats()
$ws.Cells.Item(35,6).Value = "'NO"
$ws.Cells.Item(35,6).ClearFormats()

# Row 36
$ws.Cells.Item(36,1).Value = "'3988"
$ws.Cells.Item(36,1).ClearFormats()
$ws.Cells.Item(36,2).Value = 11
$ws.Cells.Item(36,2).ClearFormats()
$ws.Cells.Item(36,3).Value = "'0"
$ws.Cells.Item(36,3).ClearFormats()
$ws.Cells.Item(36,4).Value = "'0"
$ws.Cells.Item(36,4).ClearFormats()
$ws.Cells.Item(36,5).Value = "'0.35%"
$ws.Cells.Item(36,5).ClearFormats()
$ws.Cells.Item(36,6).Value = "'NO"
$ws.Cells.Item(36,6).ClearFormats()

# Row 37
$ws.Cells.Item(37,1).Value = "'4032"
$ws.Cells.Item(37,1).ClearFormats()
$ws.Cells.Item(37,2).Value = "'"
$ws.Cells.Item(37,2).ClearFormats()
$ws.Cells.Item(37,3).Value = "'"
$ws.Cells.Item(37,3).ClearFormats()
$ws.Cells.Item(37,4).Value = "'"
$ws.Cells.Item(37,4).ClearFormats()
$ws.Cells.Item(37,5).Value = "'"
$ws.Cells.Item(37,5).ClearFormats()
$ws.Cells.Item(37,6).Value = "'NO"
$ws.Cells.Item(37,6).ClearFormats()

# Row 38
$ws.Cells.Item(38,1).Value = "'4035"
$ws.Cells.Item(38,1).ClearFormats()
$ws.Cells.Item(38,2).Value = 11
$ws.Cells.Item(38,2).ClearFormats()
$ws.Cells.Item(38,3).Value = "'"
$ws.Cells.Item(38,3).ClearFormats()
$ws.Cells.Item(38,4).Value = "'"
$ws.Cells.Item(38,4).ClearFormats()
$ws.Cells.Item(38,5).Value = "'"
$ws.Cells.Item(38,5).ClearFormats()
$ws.Cells.Item(38,6).Value = "'NO"
$ws.Cells.Item(38,6).ClearFormats()

# Row 39
$ws.Cells.Item(39,1).Value = "'4041"
$ws.Cells.Item(39,1).ClearFormats()
$ws.Cells.Item(39,2).Value = 11
$ws.Cells.Item(39,2).ClearFormats()
$ws.Cells.Item(39,3).Value = "'0"
$ws.Cells.Item(39,3).ClearFormats()
$ws.Cells.Item(39,4).Value = "'0"
$ws.Cells.Item(39,4).ClearFormats()
$ws.Cells.Item(39,5).Value = "'"
$ws.Cells.Item(39,5).ClearFormats()
$ws.Cells.Item(39,6).Value = "'YES"
$ws.Cells.Item(39,6).ClearFormats()

# Row 40
$ws.Cells.Item(40,1).Value = "'4117"
$ws.Cells.Item(40,1).ClearFormats()
$ws.Cells.Item(40,2).Value = "'"
$ws.Cells.Item(40,2).ClearFormats()
$ws.Cells.Item(40,3).Value = "'"
$ws.Cells.Item(40,3).ClearFormats()
$ws.Cells.Item(40,4).Value = "'"
$ws.Cells.Item(40,4).ClearFormats()
$ws.Cells.Item(40,5).Value = "'"
$ws.Cells.Item(40,5).ClearFormats()
$ws.Cells.Item(40,6).Value = "'NO"
$ws.Cells.Item(40,6).ClearFormats()

# Row 41
$ws.Cells.Item(41,1).Value = "'4123"
$ws.Cells.Item(41,1).ClearFormats()
$ws.Cells.Item(41,2).Value = "'"
$ws.Cells.Item(41,2).ClearFormats()
$ws.Cells.Item(41,3).Value = "'"
$ws.Cells.Item(41,3).ClearFormats()
$ws.Cells.Item(41,4).Value = "'"
$ws.Cells.Item(41,4).ClearFormats()
$ws.Cells.Item(41,5).Value = "'"
$ws.Cells.Item(41,5).ClearFormats()
$ws.Cells.Item(41,6).Value = "'NO"
$ws.Cells.Item(41,6).ClearFormats()

# Row 42
$ws.Cells.Item(42,1).Value = "'4125"
$ws.Cells.Item(42,1).ClearFormats()
$ws.Cells.Item(42,2).Value = 11
$ws.Cells.Item(42,2).ClearFormats()
$ws.Cells.Item(42,3).Value = "'0"
$ws.Cells.Item(42,3).ClearFormats()
$ws.Cells.Item(42,4).Value = "'0"
$ws.Cells.Item(42,4).ClearFormats()
$ws.Cells.Item(42,5).Value = "'0.40%"
$ws.Cells.Item(42,5).ClearFormats()
$ws.Cells.Item(42,6).Value = "'NO"
$ws.Cells.Item(42,6).ClearFormats()

# Row 43
$ws.Cells.Item(43,1).Value = "'4222"
$ws.Cells.Item(43,1).ClearFormats()
$ws.Cells.Item(43,2).Value = "'"
$ws.Cells.Item(43,2).ClearFormats()
$ws.Cells.Item(43,3).Value = "'"
$ws.Cells.Item(43,3).ClearFormats()
$ws.Cells.Item(43,4).Value = "'"
$ws.Cells.Item(43,4).ClearFormats()
$ws.Cells.Item(43,5).Value = "'"
$ws.Cells.Item(43,5).ClearFormats()
$ws.Cells.Item(43,6).Value = "'NO"
$ws.Cells.Item(43,6).ClearFormats()

# Row 44
$ws.Cells.Item(44,1).Value = "'4224"
$ws.Cells.Item(44,1).ClearFormats()
$ws.Cells.Item(44,2).Value = "'"
$ws.Cells.Item(44,2).ClearFormats()
$ws.Cells.Item(44,3).Value = "'"
$ws.Cells.Item(44,3).ClearFormats()
$ws.Cells.Item(44,4).Value = "'"
$ws.Cells.Item(44,4).ClearFormats()
$ws.Cells.Item(44,5).Value = "'"
$ws.Cells.Item(44,5).ClearFormats()
$ws.Cells.Item(44,6).Value = "'NO"
$ws.Cells.Item(44,6).ClearFormats()

# Row 45
$ws.Cells.Item(45,1).Value = "'4226"
$ws.Cells.Item(45,1).ClearFormats()
$ws.Cells.Item(45,2).Value = "'"
$ws.Cells.Item(45,2).ClearFormats()
$ws.Cells.Item(45,3).Value = "'"
$ws.Cells.Item(45,3).ClearFormats()
$ws.Cells.Item(45,4).Value = "'"
$ws.Cells.Item(45,4).ClearFormats()
$ws.Cells.Item(45,5).Value = "'"
$ws.Cells.Item(45,5).ClearFormats()
$ws.Cells.Item(45,6).Value = "'NO"
$ws.Cells.Item(45,6).ClearFormats()

# Row 46
$ws.Cells.Item(46,1).Value = "'4400"
$ws.Cells.Item(46,1).ClearFormats()
$ws.Cells.Item(46,2).Value = 11
$ws.Cells.Item(46,2).ClearFormats()
$ws.Cells.Item(46,3).Value = "'0"
$ws.Cells.Item(46,3).ClearFormats()
$ws.Cells.Item(46,4).Value = "'0"
$ws.Cells.Item(46,4).ClearFormats()
$ws.Cells.Item(46,5).Value = "'"
$ws.Cells.Item(46,5).ClearFormats()
$ws.Cells.Item(46,6).Value = "'NO"
$ws.Cells.Item(46,6).ClearFormats()

# Row 47
$ws.Cells.Item(47,1).Value = "'4415"
$ws.Cells.Item(47,1).ClearFormats()
$ws.Cells.Item(47,2).Value = "'"
$ws.Cells.Item(47,2).ClearFormats()
$ws.Cells.Item(47,3).Value = "'"
$ws.Cells.Item(47,3).ClearFormats()
$ws.Cells.Item(47,4).Value = "'"
$ws.Cells.Item(47,4).ClearFormats()
$ws.Cells.Item(47,5).Value = "'"
$ws.Cells.Item(47,5).ClearFormats()
$ws.Cells.Item(47,6).Value = "'NO"
$ws.Cells.Item(47,6).ClearFormats()

# Row 48
$ws.Cells.Item(48,1).Value = "'4421"
$ws.Cells.Item(48,1).ClearFormats()
$ws.Cells.Item(48,2).Value = 11
$ws.Cells.Item(48,2).ClearFormats()
$ws.Cells.Item(48,3).Value = "'1"
$ws.Cells.Item(48,3).ClearFormats()
$ws.Cells.Item(48,4).Value = "'1"
$ws.Cells.Item(48,4).ClearFormats()
$ws.Cells.Item(48,5).Value = "'11.79%"
$ws.Cells.Item(48,5).ClearFormats()
$ws.Cells.Item(48,6).Value = "'NO"
$ws.Cells.Item(48,6).ClearFormats()

# Row 49
$ws.Cells.Item(49,1).Value = "'4423"
$ws.Cells.Item(49,1).ClearFormats()
$ws.Cells.Item(49,2).Value = 11
$ws.Cells.Item(49,2).ClearFormats()
$ws.Cells.Item(49,3).Value = "'"
$ws.Cells.Item(49,3).ClearFormats()
$ws.Cells.Item(49,4).Value = "'"
$ws.Cells.Item(49,4).ClearFormats()
$ws.Cells.Item(49,5).Value = "'"
$ws.Cells.Item(49,5).ClearFormats()
$ws.Cells.Item(49,6).Value = "'NO"
$ws.Cells.Item(49,6).ClearFormats()

# Row 50
$ws.Cells.Item(50,1).Value = "'4429"
$ws.Cells.Item(50,1).ClearFormats()
$ws.Cells.Item(50,2).Value = "'"
$ws.Cells.Item(50,2).ClearFormats()
$ws.Cells.Item(50,3).Value = "'"
$ws.Cells.Item(50,3).ClearFormats()
$ws.Cells.Item(50,4).Value = "'"
$ws.Cells.Item(50,4).ClearFormats()
$ws.Cells.Item(50,5).Value = "'"
$ws.Cells.Item(50,5).ClearFormats()
$ws.Cells.Item(50,6).Value = "'NO"
$ws.Cells.Item(50,6).ClearFormats()

# Row 51
$ws.Cells.Item(51,1).Value = "'4430"
$ws.Cells.Item(51,1).ClearFormats()
$ws.Cells.Item(51,2).Value = "'"
$ws.Cells.Item(51,2).ClearFormats()
$ws.Cells.Item(51,3).Value = "'"
$ws.Cells.Item(51,3).ClearFormats()
$ws.Cells.Item(51,4).Value = "'"
$ws.Cells.Item(51,4).ClearFormats()
$ws.Cells.Item(51,5).Value = "'"
$ws.Cells.Item(51,5).ClearFormats()
$ws.Cells.Item(51,6).Value = "'NO"
$ws.Cells.Item(51,6).ClearFormats()

# Row 52
$ws.Cells.Item(52,1).Value = "'4431"
$ws.Cells.Item(52,1).ClearFormats()
$ws.Cells.Item(52,2).Value = "'"
$ws.Cells.Item(52,2).ClearFormats()
$ws.Cells.Item(52,3).Value = "'"
$ws.Cells.Item(52,3).ClearFormats()
$ws.Cells.Item(52,4).Value = "'"
$ws.Cells.Item(52,4).ClearFormats()
$ws.Cells.Item(52,5).Value = "'"
$ws.Cells.Item(52,5).ClearFormats()
$ws.Cells.Item(52,6).Value = "'"
$ws.Cells.Item(52,6).ClearFormats()

# Row 53
$ws.Cells.Item(53,1).Value = "'4435"
$ws.Cells.Item(53,1).ClearFormats()
$ws.Cells.Item(53,2).Value = "'"
$ws.Cells.Item(53,2).ClearFormats()
$ws.Cells.Item(53,3).Value = "'"
$ws.Cells.Item(53,3).ClearFormats()
$ws.Cells.Item(53,4).Value = "'"
$ws.Cells.Item(53,4).ClearFormats()
$ws.Cells.Item(53,5).Value = "'"
$ws.Cells.Item(53,5).ClearFormats()
$ws.Cells.Item(53,6).Value = "'"
$ws.Cells.Item(53,6).ClearFormats()

# Row 54
$ws.Cells.Item(54,1).Value = "'4436"
$ws.Cells.Item(54,1).ClearFormats()
$ws.Cells.Item(54,2).Value = "'"
$ws.Cells.Item(54,2).ClearFormats()
$ws.Cells.Item(54,3).Value = "'"
$ws.Cells.Item(54,3).ClearFormats()
$ws.Cells.Item(54,4).Value = "'"
$ws.Cells.Item(54,4).ClearFormats()
$ws.Cells.Item(54,5).Value = "'"
$ws.Cells.Item(54,5).ClearFormats()
$ws.Cells.Item(54,6).Value = "'"
$ws.Cells.Item(54,6).ClearFormats()

# Row 55
$ws.Cells.Item(55,1).Value = "'4437"
$ws.Cells.Item(55,1).ClearFormats()
$ws.Cells.Item(55,2).Value = "'"
$ws.Cells.Item(55,2).ClearFormats()
$ws.Cells.Item(55,3).Value = "'"
$ws.Cells.Item(55,3).ClearFormats()
$ws.Cells.Item(55,4).Value = "'"
$ws.Cells.Item(55,4).ClearFormats()
$ws.Cells.Item(55,5).Value = "'"
$ws.Cells.Item(55,5).ClearFormats()
$ws.Cells.Item(55,6).Value = "'"
$ws.Cells.Item(55,6).ClearFormats()

# Row 56
$ws.Cells.Item(56,1).Value = "'4483"
$ws.Cells.Item(56,1).ClearFormats()
$ws.Cells.Item(56,2).Value = "'"
$ws.Cells.Item(56,2).ClearFormats()
$ws.Cells.Item(56,3).Value = "'"
$ws.Cells.Item(56,3).ClearFormats()
$ws.Cells.Item(56,4).Value = "'"
$ws.Cells.Item(56,4).ClearFormats()
$ws.Cells.Item(56,5).Value = "'"
$ws.Cells.Item(56,5).ClearFormats()
$ws.Cells.Item(56,6).Value = "'"
$ws.Cells.Item(56,6).ClearFormats()

# Row 57
$ws.Cells.Item(57,1).Value = "'4486"
$ws.Cells.Item(57,1).ClearFormats()
$ws.Cells.Item(57,2).Value = "'"
$ws.Cells.Item(57,2).ClearFormats()
$ws.Cells.Item(57,3).Value = "'"
$ws.Cells.Item(57,3).ClearFormats()
$ws.Cells.Item(57,4).Value = "'"
$ws.Cells.Item(57,4).ClearFormats()
$ws.Cells.Item(57,5).Value = "'"
$ws.Cells.Item(57,5).ClearFormats()
$ws.Cells.Item(57,6).Value = "'"
$ws.Cells.Item(57,6).ClearFormats()

# Row 58
$ws.Cells.Item(58,1).Value = "'4594"
$ws.Cells.Item(58,1).ClearFormats()
$ws.Cells.Item(58,2).Value = "'"
$ws.Cells.Item(58,2).ClearFormats()
$ws.Cells.Item(58,3).Value = "'"
$ws.Cells.Item(58,3).ClearFormats()
$ws.Cells.Item(58,4).Value = "'"
$ws.Cells.Item(58,4).ClearFormats()
$ws.Cells.Item(58,5).Value = "'"
$ws.Cells.Item(58,5).ClearFormats()
$ws.Cells.Item(58,6).Value = "'"
$ws.Cells.Item(58,6).ClearFormats()

# Row 59
$ws.Cells.Item(59,1).Value = "'4597"
$ws.Cells.Item(59,1).ClearFormats()
$ws.Cells.Item(59,2).Value = "'"
$ws.Cells.Item(59,2).ClearFormats()
$ws.Cells.Item(59,3).Value = "'"
$ws.Cells.Item(59,3).ClearFormats()
$ws.Cells.Item(59,4).Value = "'"
$ws.Cells.Item(59,4).ClearFormats()
$ws.Cells.Item(59,5).Value = "'"
$ws.Cells.Item(59,5).ClearFormats()
$ws.Cells.Item(59,6).Value = "'"
$ws.Cells.Item(59,6).ClearFormats()

# Row 60
$ws.Cells.Item(60,1).Value = "'4600"
$ws.Cells.Item(60,1).ClearFormats()
$ws.Cells.Item(60,2).Value = "'"
$ws.Cells.Item(60,2).ClearFormats()
$ws.Cells.Item(60,3).Value = "'"
$ws.Cells.Item(60,3).ClearFormats()
$ws.Cells.Item(60,4).Value = "'"
$ws.Cells.Item(60,4).ClearFormats()
$ws.Cells.Item(60,5).Value = "'"
$ws.Cells.Item(60,5).ClearFormats()
$ws.Cells.Item(60,6).Value = "'"
$ws.Cells.Item(60,6).ClearFormats()

# Row 61
$ws.Cells.Item(61,1).Value = "'4601"
$ws.Cells.Item(61,1).ClearFormats()
$ws.Cells.Item(61,2).Value = "'"
$ws.Cells.Item(61,2).ClearFormats()
$ws.Cells.Item(61,3).Value = "'"
$ws.Cells.Item(61,3).ClearFormats()
$ws.Cells.Item(61,4).Value = "'"
$ws.Cells.Item(61,4).ClearFormats()
$ws.Cells.Item(61,5).Value = "'"
$ws.Cells.Item(61,5).ClearFormats()
$ws.Cells.Item(61,6).Value = "'"
$ws.Cells.Item(61,6).ClearFormats()

# Row 62
$ws.Cells.Item(62,1).Value = "'4603"
$ws.Cells.Item(62,1).ClearFormats()
$ws.Cells.Item(62,2).Value = "'"
$ws.Cells.Item(62,2).ClearFormats()
$ws.Cells.Item(62,3).Value = "'"
$ws.Cells.Item(62,3).ClearFormats()
$ws.Cells.Item(62,4).Value = "'"
$ws.Cells.Item(62,4).ClearFormats()
$ws.Cells.Item(62,5).Value = "'"
$ws.Cells.Item(62,5).ClearFormats()
$ws.Cells.Item(62,6).Value = "'"
$ws.Cells.Item(62,6).ClearFormats()

# Row 63
$ws.Cells.Item(63,1).Value = "'4644"
$ws.Cells.Item(63,1).ClearFormats()
$ws.Cells.Item(63,2).Value = "'"
$ws.Cells.Item(63,2).ClearFormats()
$ws.Cells.Item(63,3).Value = "'"
$ws.Cells.Item(63,3).ClearFormats()
$ws.Cells.Item(63,4).Value = "'"
$ws.Cells.Item(63,4).ClearFormats()
$ws.Cells.Item(63,5).Value = "'"
$ws.Cells.Item(63,5).ClearFormats()
$ws.Cells.Item(63,6).Value = "'"
$ws.Cells.Item(63,6).ClearFormats()

# Row 64
$ws.Cells.Item(64,1).Value = "'4645"
$ws.Cells.Item(64,1).ClearFormats()
$ws.Cells.Item(64,2).Value = "'"
$ws.Cells.Item(64,2).ClearFormats()
$ws.Cells.Item(64,3).Value = "'"
$ws.Cells.Item(64,3).ClearFormats()
$ws.Cells.Item(64,4).Value = "'"
$ws.Cells.Item(64,4).ClearFormats()
$ws.Cells.Item(64,5).Value = "'"
$ws.Cells.Item(64,5).ClearFormats()
$ws.Cells.Item(64,6).Value = "'"
$ws.Cells.Item(64,6).ClearFormats()

# Row 65
$ws.Cells.Item(65,1).Value = "'4646"
$ws.Cells.Item(65,1).ClearFormats()
$ws.Cells.Item(65,2).Value = "'"
$ws.Cells.Item(65,2).ClearFormats()
$ws.Cells.Item(65,3).Value = "'"
$ws.Cells.Item(65,3).ClearFormats()
$ws.Cells.Item(65,4).Value = "'"
$ws.Cells.Item(65,4).ClearFormats()
$ws.Cells.Item(65,5).Value = "'"
$ws.Cells.Item(65,5).ClearFormats()
$ws.Cells.Item(65,6).Value = "'"
$ws.Cells.Item(65,6).ClearFormats()

# Row 66
$ws.Cells.Item(66,1).Value = "'4647"
$ws.Cells.Item(66,1).ClearFormats()
$ws.Cells.Item(66,2).Value = "'"
$ws.Cells.Item(66,2).ClearFormats()
$ws.Cells.Item(66,3).Value = "'"
$ws.Cells.Item(66,3).ClearFormats()
$ws.Cells.Item(66,4).Value = "'"
$ws.Cells.Item(66,4).ClearFormats()
$ws.Cells.Item(66,5).Value = "'"
$ws.Cells.Item(66,5).ClearFormats()
$ws.Cells.Item(66,6).Value = "'"
$ws.Cells.Item(66,6).ClearFormats()

# Row 67
$ws.Cells.Item(67,1).Value = "'4648"
$ws.Cells.Item(67,1).ClearFormats()
$ws.Cells.Item(67,2).Value = "'"
$ws.Cells.Item(67,2).ClearFormats()
$ws.Cells.Item(67,3).Value = "'"
$ws.Cells.Item(67,3).ClearFormats()
$ws.Cells.Item(67,4).Value = "'"
$ws.Cells.Item(67,4).ClearFormats()
$ws.Cells.Item(67,5).Value = "'"
$ws.Cells.Item(67,5).ClearFormats()
$ws.Cells.Item(67,6).Value = "'"
$ws.Cells.Item(67,6).ClearFormats()

# Row 68
$ws.Cells.Item(68,1).Value = "'4649"
$ws.Cells.Item(68,1).ClearFormats()
$ws.Cells.Item(68,2).Value = "'"
$ws.Cells.Item(68,2).ClearFormats()
$ws.Cells.Item(68,3).Value = "'"
$ws.Cells.Item(68,3).ClearFormats()
$ws.Cells.Item(68,4).Value = "'"
$ws.Cells.Item(68,4).ClearFormats()
$ws.Cells.Item(68,5).Value = "'"
$ws.Cells.Item(68,5).ClearFormats()
$ws.Cells.Item(68,6).Value = "'"
$ws.Cells.Item(68,6).ClearFormats()

# Row 69
$ws.Cells.Item(69,1).Value = "'4663"
$ws.Cells.Item(69,1).ClearFormats()
$ws.Cells.Item(69,2).Value = "'"
$ws.Cells.Item(69,2).ClearFormats()
$ws.Cells.Item(69,3).Value = "'"
$ws.Cells.Item(69,3).ClearFormats()
$ws.Cells.Item(69,4).Value = "'"
$ws.Cells.Item(69,4).ClearFormats()
$ws.Cells.Item(69,5).Value = "'"
$ws.Cells.Item(69,5).ClearFormats()
$ws.Cells.Item(69,6).Value = "'"
$ws.Cells.Item(69,6).ClearFormats()

# Row 70
$ws.Cells.Item(70,1).Value = "'4666"
$ws.Cells.Item(70,1).ClearFormats()
$ws.Cells.Item(70,2).Value = "'"
$ws.Cells.Item(70,2).ClearFormats()
$ws.Cells.Item(70,3).Value = "'"
$ws.Cells.Item(70,3).ClearFormats()
$ws.Cells.Item(70,4).Value = "'"
$ws.Cells.Item(70,4).ClearFormats()
$ws.Cells.Item(70,5).Value = "'"
$ws.Cells.Item(70,5).ClearFormats()
$ws.Cells.Item(70,6).Value = "'"
$ws.Cells.Item(70,6).ClearFormats()
